# Insert a new row at position 77, shifting existing rows 77..133 down to 78..134
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new data record
$ws.Range("A77").Value = 5
$ws.Range("B77").Value = "Macroferia Regional de Talca"
$ws.Range("C77").Value = "Maule"
$ws.Range("D77").Value = 44741
$ws.Range("E77").Value = 7
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108002
$ws.Range("J77").Value = "Mango"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 240
$ws.Range("N77").Value = 8000
$ws.Range("O77").Value = 8000
$ws.Range("P77").Value = 8000
$ws.Range("Q77").Value = "$/bandeja 4 kilos"
$ws.Range("R77").Value = "Brasil"
$ws.Range("S77").Value = 2000
$ws.Range("T77").Value = 4
